$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")

# --- Fix the misspelled "Checkpoitn 2224" entry at A142 ---
# Clearing first removes the old shared-string entry entirely (it was only
# referenced once), which re-numbers the shared-string table exactly the way
# the target file does (Checkpoint 2586 / Touch button / Enter 8-5 each shift
# down by one). Re-entering the corrected spelling then appends a brand new
# shared-string entry at the end of the table.
$ws.Range("A142").Value = ""
$ws.Range("A142").Value = "Checkpoint 2224"

# --- Append the new rows (149-167) capturing the rest of 8-5 and the start of 8-6 ---
# The order the A-column labels are written in matters: it reproduces the
# exact order new entries were appended to the shared-strings table in the
# source commit (notably "Jump block" was registered before "Jump Pipe").

$ws.Range("A149").Value = "Checkpoint 271"
$ws.Range("B149").Value = 50797
$ws.Range("C149").Value = 59960
$ws.Range("D149").Formula = "=IF(B149 >  0,C149-B149, 0)"

$ws.Range("A150").Value = "Checkpoint 596/595"
$ws.Range("B150").Value = 50904
$ws.Range("C150").Value = 60068
$ws.Range("D150").Formula = "=IF(B150 >  0,C150-B150, 0)"

$ws.Range("A151").Value = "Checkoint 872/870"
$ws.Range("B151").Value = 50995
$ws.Range("C151").Value = 60159
$ws.Range("D151").Formula = "=IF(B151 >  0,C151-B151, 0)"

$ws.Range("A152").Value = "Checkpoint 1293/1291"
$ws.Range("B152").Value = 51135
$ws.Range("C152").Value = 60299
$ws.Range("D152").Formula = "=IF(B152 >  0,C152-B152, 0)"

$ws.Range("A153").Value = "Checkpoint 1945"
$ws.Range("B153").Value = 51352
$ws.Range("C153").Value = 60517
$ws.Range("D153").Formula = "=IF(B153 >  0,C153-B153, 0)"

$ws.Range("A154").Value = "Checkpoint 2354"
$ws.Range("B154").Value = 51487
$ws.Range("C154").Value = 60652
$ws.Range("D154").Formula = "=IF(B154 >  0,C154-B154, 0)"

$ws.Range("A155").Value = "Checkpoitn 2941"
$ws.Range("B155").Value = 51681
$ws.Range("C155").Value = 60846
$ws.Range("D155").Formula = "=IF(B155 >  0,C155-B155, 0)"

$ws.Range("A156").Value = "Enter Pipe"
$ws.Range("B156").Value = 52080
$ws.Range("C156").Value = 61248
$ws.Range("D156").Formula = "=IF(B156 >  0,C156-B156, 0)"

$ws.Range("A157").Value = "Get Flag"
$ws.Range("B157").Value = 52250
$ws.Range("C157").Value = 61427
$ws.Range("D157").Formula = "=IF(B157 >  0,C157-B157, 0)"

$ws.Range("A158").Value = "End Level"
$ws.Range("B158").Value = 52768
$ws.Range("C158").Value = 61945
$ws.Range("D158").Formula = "=IF(B158 >  0,C158-B158, 0)"

$ws.Range("A159").Value = "Enter 8-6"
$ws.Range("B159").Value = 53104
$ws.Range("C159").Value = 62610
$ws.Range("D159").Formula = "=IF(B159 >  0,C159-B159, 0)"

$ws.Range("A160").Value = "1st Move"
$ws.Range("B160").Value = 53332
$ws.Range("C160").Value = 62861
$ws.Range("D160").Formula = "=IF(B160 >  0,C160-B160, 0)"

$ws.Range("A161").Value = "Land 1st Orange plat"
$ws.Range("B161").Value = 53477
$ws.Range("C161").Value = 63021
$ws.Range("D161").Formula = "=IF(B161 >  0,C161-B161, 0)"

# "Jump block" (row 166) is registered in the shared-strings table before
# "Jump Pipe" (row 162) in the source commit, so write it first.
$ws.Range("A166").Value = "Jump block"

$ws.Range("A162").Value = "Jump Pipe"
$ws.Range("B162").Value = 53562
$ws.Range("C162").Value = 63110
$ws.Range("D162").Formula = "=IF(B162 >  0,C162-B162, 0)"

$ws.Range("B163").Value = 53643
$ws.Range("C163").Value = 63206
$ws.Range("D163").Formula = "=IF(B163 >  0,C163-B163, 0)"

$ws.Range("A164").Value = "Checkpoint"
$ws.Range("B164").Value = 53813
$ws.Range("C164").Value = 63376
$ws.Range("D164").Formula = "=IF(B164 >  0,C164-B164, 0)"

$ws.Range("A165").Value = "Push on spring"
$ws.Range("B165").Value = 53851
$ws.Range("C165").Value = 63425
$ws.Range("D165").Formula = "=IF(B165 >  0,C165-B165, 0)"

# Row 166's B/C/D (A166 was already set above, ahead of row 162, to match
# shared-string append order).
$ws.Range("B166").Value = 53952
$ws.Range("C166").Value = 63540
$ws.Range("D166").Formula = "=IF(B166 >  0,C166-B166, 0)"

$ws.Range("A167").Value = "Spring off 2nd spring (sparks)"
$ws.Range("B167").Value = 54007
$ws.Range("C167").Value = 63597
$ws.Range("D167").Formula = "=IF(B167 >  0,C167-B167, 0)"

# --- Update the frozen-pane scroll position / active selection to match ---
$ws.Range("B168").Select()
